# Integration plan: correct the "who's working on which step" markers
# (T = currently being worked on, X = already done, S = scheduled/not started
# yet) and drop the now-obsolete final row (the plan finishes one step
# earlier than before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Step 1) ---
$ws.Range("D2").Value = "T"
$ws.Range("E2").Value = "X"
$ws.Range("F2").Value = "S"
$ws.Range("G2").Value = "S"
$ws.Range("H2").Value = "S"

# --- Row 3 (Step 2) ---
$ws.Range("D3").Value = "T"
$ws.Range("E3").Value = "X"
$ws.Range("F3").Value = "X"
$ws.Range("G3").Value = "S"
$ws.Range("H3").Value = "S"

# --- Row 4 (Step 3) ---
$ws.Range("D4").Value = "T"
$ws.Range("E4").Value = "X"
$ws.Range("F4").Value = "X"
$ws.Range("G4").Value = "X"
$ws.Range("H4").Value = "S"

# --- Row 5 (Step 4) ---
$ws.Range("D5").Value = "T"
$ws.Range("E5").Value = "X"
$ws.Range("F5").Value = "X"
$ws.Range("G5").Value = "X"
$ws.Range("H5").Value = "X"

# --- Row 6 (Step 5) ---
$ws.Range("C6").Value = "T"
$ws.Range("D6").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("F6").Value = "X"
$ws.Range("G6").Value = "X"
$ws.Range("H6").Value = "X"

# --- Row 7 (Step 6) ---
$ws.Range("B7").Value = "T"
$ws.Range("C7").Value = "X"
$ws.Range("D7").Value = "X"
$ws.Range("E7").Value = "X"
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "X"
$ws.Range("H7").Value = "X"

# Plan now ends at step 6 - remove the old final "step 7" row (this also
# shrinks the dimension, the table range and the autofilter automatically).
$ws.Rows(8).Delete()

# Column widths shift slightly (bestFit) now the header row is the widest
# content in every column; nudge them to the new best-fit pixel widths.
$ws.Columns(1).ColumnWidth = 6.66666666666667
$ws.Columns(2).ColumnWidth = 21.6666666666667
$ws.Columns(3).ColumnWidth = 10
$ws.Columns(4).ColumnWidth = 14.3333333333333
$ws.Columns(5).ColumnWidth = 17.8333333333333
$ws.Columns(6).ColumnWidth = 19.6666666666667
$ws.Columns(7).ColumnWidth = 24.6666666666667

# View moved on: scrolled right one column and zoomed out, with the
# selection now resting further down the (now shorter) sheet.
$excel.ActiveWindow.Zoom = 86
$ws.Range("D14").Select()
